$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column M ("In Advance"), shifting
# In Advance / Late / (blank heading) / Outstanding one column to the right.
$ws.Columns("M").Insert()

# Give the newly inserted column a custom width similar to its neighbours.
$ws.Columns("M").ColumnWidth = 8.33

# Update the "Due" total for row 3 (was 10145.16, now split out to 145.16).
$ws.Range("K3").Value = 145.16

# The (now shifted) "In Advance" column for row 3 gets the principal amount.
$ws.Range("N3").Value = 10000

# Make "Repayment schedule" the active/selected sheet and update its selection.
$ws.Activate()
$ws.Range("I16").Select()

$win = $excel.ActiveWindow
$win.Zoom = 100
